# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets.
#
# - Insert a new "Player Info" worksheet before the existing "ODI Batting"
#   sheet and populate it with the player's basic details.
# - On the "ODI Batting" sheet, simplify the MATCH_CARD_LINK column (full
#   scorecard URLs) down to a MATCH_CODE column (just the numeric code).

$wb = $excel.ActiveWorkbook

$odiBatting = $wb.Worksheets.Item("ODI Batting")

# Insert the new sheet immediately before "ODI Batting" so the tab order
# becomes: Player Info, ODI Batting.
$playerInfo = $wb.Worksheets.Add($odiBatting)
$playerInfo.Name = "Player Info"

# Header row.
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Data row for player 6470 (Harry Brook).
$playerInfo.Range("A2").Value = "6470"
$playerInfo.Range("B2").Value = "Harry Cherrington Brook"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Medium"

$playerInfo.Range("A1").Select()

# Update the ODI Batting sheet: rename MATCH_CARD_LINK -> MATCH_CODE and
# replace the full scorecard URLs with just the bare match code.
# (Re-fetch the sheet by name: the handle captured before Worksheets.Add()
# can no longer be trusted to point at the right sheet afterwards.)
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("D1").Value = "MATCH_CODE"
$odiBatting.Range("D2").Value = "4698"
$odiBatting.Range("D3").Value = "4699"
$odiBatting.Range("D4").Value = "4700"
